$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats
$xlPasteFormats = -4122

# ------------------------------------------------------------------
# 1. Fix the dates on rows 31-33: they were all entered as 45600
#    (11/4/2024) but the first two actually happened on 45599
#    (11/3/2024).
# ------------------------------------------------------------------
$ws.Range("A31").Value = 45599
$ws.Range("A32").Value = 45599
$ws.Range("A33").Value = 45599

# ------------------------------------------------------------------
# 2. Add the new time-sheet entry (row 34) for 11/5/2024 - the shift
#    spent creating the new friction pads.
# ------------------------------------------------------------------

# Copy formatting from existing rows that already have the desired
# number formats so we re-use the same style indices instead of
# minting new (duplicate) styles.
$ws.Range("A2").Copy()
$ws.Range("A34").PasteSpecial($xlPasteFormats)

$ws.Range("B17").Copy()
$ws.Range("B34").PasteSpecial($xlPasteFormats)

$ws.Range("C2").Copy()
$ws.Range("C34").PasteSpecial($xlPasteFormats)

$ws.Range("D2").Copy()
$ws.Range("D34").PasteSpecial($xlPasteFormats)

$ws.Range("E2").Copy()
$ws.Range("E34").PasteSpecial($xlPasteFormats)

$ws.Range("F2").Copy()
$ws.Range("F34").PasteSpecial($xlPasteFormats)

$ws.Range("G2").Copy()
$ws.Range("G34").PasteSpecial($xlPasteFormats)

$ws.Range("A34").Value = 45601
$ws.Range("B34").Value = 0.89583333333333337
$ws.Range("C34").Value = 0.94444444444444442
$ws.Range("D34").Formula = "=C34-B34"
$ws.Range("E34").Formula = "=D34*1440"
$ws.Range("F34").Formula = "=E34/60"
$ws.Range("G34").Formula = "=F34*22.5"

# ------------------------------------------------------------------
# 3. Weekly roll-up for the new (eighth) week -- mirrors the pattern
#    already used in M3:M7 / N3:N7 for the previous weeks.
# ------------------------------------------------------------------
$ws.Range("M7").Copy()
$ws.Range("M8").PasteSpecial($xlPasteFormats)

$ws.Range("N7").Copy()
$ws.Range("N8").PasteSpecial($xlPasteFormats)

$ws.Range("M8").Formula = "=SUM(D34)"
$ws.Range("N8").Formula = "=SUM(G34)"

# ------------------------------------------------------------------
# 4. Move the "DON'T FORGET TO SAVE (local)" reminder down from row
#    36 to row 46 to leave room for future entries.
# ------------------------------------------------------------------
$reminder = $ws.Range("B36").Value2
$ws.Rows("36:36").Delete()

$ws.Range("B1").Copy()
$ws.Range("B46").PasteSpecial($xlPasteFormats)
$ws.Range("B46").Value = $reminder

# ------------------------------------------------------------------
# 5. Tidy up the view: drop the frozen/scroll anchor and move the
#    active selection to where we were last working.
# ------------------------------------------------------------------
$ws.Range("K14").Select()

$wb.Save()
